# Auto-generated edit script: updates market-board derived columns (H-N)
# on specific rows across all 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to refresh currentAveragePrice* / LevePrice* / LeveProfit* figures.

$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
# row 4
$ws.Range("H4").Value = 91.25
$ws.Range("I4").Value = 75.71429000000001
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 75.71429000000001
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = 38.28570999999999
$ws.Range("N4").Value = -428
# row 98
$ws.Range("H98").Value = 778.05
$ws.Range("I98").Value = 778.05
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 778.05
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 719.95
$ws.Range("N98").Value = $null
# row 100
$ws.Range("H100").Value = 5948.2085
$ws.Range("I100").Value = 3846.6667
$ws.Range("J100").Value = 6648.722
$ws.Range("K100").Value = 3846.6667
$ws.Range("L100").Value = 6648.722
$ws.Range("M100").Value = -3305.6667
$ws.Range("N100").Value = -7730.722
# row 111
$ws.Range("H111").Value = 657
$ws.Range("I111").Value = 515.8
$ws.Range("J111").Value = 858.7143
$ws.Range("K111").Value = 1547.4
$ws.Range("L111").Value = 2576.1429
$ws.Range("M111").Value = 1519.6
$ws.Range("N111").Value = -8710.142899999999
# row 122
$ws.Range("H122").Value = 778.05
$ws.Range("I122").Value = 778.05
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2334.15
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 115.8500000000004
$ws.Range("N122").Value = $null
# row 132
$ws.Range("H132").Value = 2060.908
$ws.Range("I132").Value = 1209.5883
$ws.Range("K132").Value = 3628.7649
$ws.Range("M132").Value = -1098.7649
# row 135
$ws.Range("H135").Value = 469.4737
$ws.Range("I135").Value = 300.66666
$ws.Range("J135").Value = 1102.5
$ws.Range("K135").Value = 2705.99994
$ws.Range("L135").Value = 9922.5
$ws.Range("M135").Value = -170.9999399999997
$ws.Range("N135").Value = -14992.5
# row 137
$ws.Range("H137").Value = 2383.9524
$ws.Range("I137").Value = 2114.3674
$ws.Range("K137").Value = 6343.1022
$ws.Range("M137").Value = -3793.1022

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
# row 32
$ws.Range("H32").Value = 3363.58
$ws.Range("I32").Value = 2001.3489
$ws.Range("J32").Value = 11731.571
$ws.Range("K32").Value = 2001.3489
$ws.Range("L32").Value = 11731.571
$ws.Range("M32").Value = -1714.3489
$ws.Range("N32").Value = -12305.571
# row 74
$ws.Range("H74").Value = 156321.34
$ws.Range("I74").Value = 223011.47
$ws.Range("J74").Value = 40896.117
$ws.Range("K74").Value = 223011.47
$ws.Range("L74").Value = 40896.117
$ws.Range("M74").Value = -222137.47
$ws.Range("N74").Value = -42644.117
# row 77
$ws.Range("H77").Value = 156321.34
$ws.Range("I77").Value = 223011.47
$ws.Range("J77").Value = 40896.117
$ws.Range("K77").Value = 1115057.35
$ws.Range("L77").Value = 204480.585
$ws.Range("M77").Value = -1110689.35
$ws.Range("N77").Value = -213216.585
# row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
# row 110
$ws.Range("H110").Value = 1092.2307
$ws.Range("I110").Value = 1128.4286
$ws.Range("J110").Value = 1050
$ws.Range("K110").Value = 1128.4286
$ws.Range("L110").Value = 1050
$ws.Range("M110").Value = 916.5714
$ws.Range("N110").Value = -5140

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
# row 105
$ws.Range("H105").Value = 1815.9546
$ws.Range("I105").Value = 1664.7059
$ws.Range("K105").Value = 1664.7059
$ws.Range("M105").Value = 82.29410000000007
# row 107
$ws.Range("H107").Value = 2180.5264
$ws.Range("I107").Value = 1986.5
$ws.Range("J107").Value = 2723.8
$ws.Range("K107").Value = 1986.5
$ws.Range("L107").Value = 2723.8
$ws.Range("M107").Value = -66.5
$ws.Range("N107").Value = -6563.8
# row 112
$ws.Range("H112").Value = 39950
$ws.Range("J112").Value = 39950
$ws.Range("L112").Value = 39950
$ws.Range("N112").Value = -42904

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
# row 31
$ws.Range("H31").Value = 2647.8254
$ws.Range("I31").Value = 1579.0834
$ws.Range("J31").Value = 4072.8147
$ws.Range("K31").Value = 1579.0834
$ws.Range("L31").Value = 4072.8147
$ws.Range("M31").Value = -1284.0834
$ws.Range("N31").Value = -4662.8147
# row 34
$ws.Range("H34").Value = 2647.8254
$ws.Range("I34").Value = 1579.0834
$ws.Range("J34").Value = 4072.8147
$ws.Range("K34").Value = 1579.0834
$ws.Range("L34").Value = 4072.8147
$ws.Range("M34").Value = -1377.0834
$ws.Range("N34").Value = -4476.8147
# row 58
$ws.Range("H58").Value = 2322.25
$ws.Range("I58").Value = 2484.2454
$ws.Range("J58").Value = 1948.9565
$ws.Range("K58").Value = 2484.2454
$ws.Range("L58").Value = 1948.9565
$ws.Range("M58").Value = -2281.2454
$ws.Range("N58").Value = -2354.9565
# row 105
$ws.Range("H105").Value = 619.3333
$ws.Range("I105").Value = 592.1429000000001
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 592.1429000000001
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1154.8571
$ws.Range("N105").Value = -4494
# row 107
$ws.Range("H107").Value = 288.1579
$ws.Range("I107").Value = 253.875
$ws.Range("J107").Value = 471
$ws.Range("K107").Value = 253.875
$ws.Range("L107").Value = 471
$ws.Range("M107").Value = 1666.125
$ws.Range("N107").Value = -4311
# row 132
$ws.Range("H132").Value = 1657.7778
$ws.Range("I132").Value = 912.5
$ws.Range("J132").Value = 3961.3635
$ws.Range("K132").Value = 2737.5
$ws.Range("L132").Value = 11884.0905
$ws.Range("M132").Value = -207.5
$ws.Range("N132").Value = -16944.0905
# row 134
$ws.Range("H134").Value = 1684.3693
$ws.Range("I134").Value = 1141.3889
$ws.Range("J134").Value = 2358.4138
$ws.Range("K134").Value = 3424.1667
$ws.Range("L134").Value = 7075.241399999999
$ws.Range("M134").Value = -889.1666999999998
$ws.Range("N134").Value = -12145.2414
# row 136
$ws.Range("H136").Value = 2322.25
$ws.Range("I136").Value = 2484.2454
$ws.Range("J136").Value = 1948.9565
$ws.Range("K136").Value = 7452.736199999999
$ws.Range("L136").Value = 5846.8695
$ws.Range("M136").Value = -4902.736199999999
$ws.Range("N136").Value = -10946.8695

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
# row 5
$ws.Range("H5").Value = 523.67645
$ws.Range("I5").Value = 443.5
$ws.Range("J5").Value = 1125
$ws.Range("K5").Value = 1330.5
$ws.Range("L5").Value = 3375
$ws.Range("M5").Value = -1218.5
$ws.Range("N5").Value = -3599
# row 122
$ws.Range("H122").Value = 922.61536
$ws.Range("J122").Value = 1117.0588
$ws.Range("L122").Value = 10053.5292
$ws.Range("N122").Value = -14953.5292
# row 131
$ws.Range("H131").Value = 1883.3116
$ws.Range("J131").Value = 1731.2142
$ws.Range("L131").Value = 5193.642599999999
$ws.Range("N131").Value = -15273.6426
# row 135
$ws.Range("H135").Value = 523.67645
$ws.Range("I135").Value = 443.5
$ws.Range("J135").Value = 1125
$ws.Range("K135").Value = 3991.5
$ws.Range("L135").Value = 10125
$ws.Range("M135").Value = -1456.5
$ws.Range("N135").Value = -15195

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
# row 122
$ws.Range("H122").Value = 1181
$ws.Range("I122").Value = 1115.4
$ws.Range("J122").Value = 1399.6666
$ws.Range("K122").Value = 3346.2
$ws.Range("L122").Value = 4198.9998
$ws.Range("M122").Value = -896.2000000000003
$ws.Range("N122").Value = -9098.9998

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
# row 22
$ws.Range("H22").Value = 1178.2667
$ws.Range("I22").Value = 957.2857
$ws.Range("K22").Value = 957.2857
$ws.Range("M22").Value = -662.2857
# row 27
$ws.Range("H27").Value = 1178.2667
$ws.Range("I27").Value = 957.2857
$ws.Range("K27").Value = 957.2857
$ws.Range("M27").Value = -850.2857
# row 132
$ws.Range("H132").Value = 5104.942
$ws.Range("I132").Value = 1789.4375
$ws.Range("J132").Value = 12683.238
$ws.Range("K132").Value = 5368.3125
$ws.Range("L132").Value = 38049.714
$ws.Range("M132").Value = -2838.3125
$ws.Range("N132").Value = -43109.714
# row 136
$ws.Range("H136").Value = 2823.122
$ws.Range("I136").Value = 1565.0154
$ws.Range("J136").Value = 7633.5293
$ws.Range("K136").Value = 4695.0462
$ws.Range("L136").Value = 22900.5879
$ws.Range("M136").Value = -2145.0462
$ws.Range("N136").Value = -28000.5879

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
# row 132
$ws.Range("H132").Value = 1950.0178
$ws.Range("I132").Value = 1496.375
$ws.Range("J132").Value = 2290.25
$ws.Range("K132").Value = 4489.125
$ws.Range("L132").Value = 6870.75
$ws.Range("M132").Value = -1959.125
$ws.Range("N132").Value = -11930.75
# row 136
$ws.Range("H136").Value = 1407.79
$ws.Range("I136").Value = 1112.2239
$ws.Range("J136").Value = 2007.8788
$ws.Range("K136").Value = 3336.6717
$ws.Range("L136").Value = 6023.636399999999
$ws.Range("M136").Value = -786.6716999999999
$ws.Range("N136").Value = -11123.6364

